# Block.xlsx update
# Adds four new block entries (ids 201-204) to the end of the "Block" sheet:
#   row 204 -> id 201, EA 23.214, 熔岩方块        / magma block      / 溶岩のブロック
#   row 205 -> id 202, EA 23.238, 铁丝网栅栏       / wire mesh fence  / 金網のフェンス
#   row 206 -> id 203, EA 23.238, 铁丝网栅栏       / wire mesh fence  / 金網のフェンス
#   row 207 -> id 204, EA 23.238, 石头栅栏         / stone fence      / 石のフェンス
#
# The leading "'" on the id values forces Excel to store them as text
# (matching the existing id column, which is text-typed) instead of numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(204, 1).Value = "'201"
$ws.Cells.Item(204, 2).Value = "EA 23.214"
$ws.Cells.Item(204, 3).Value = "熔岩方块"
$ws.Cells.Item(204, 4).Value = "magma block"
$ws.Cells.Item(204, 5).Value = "溶岩のブロック"

$ws.Cells.Item(205, 1).Value = "'202"
$ws.Cells.Item(205, 2).Value = "EA 23.238"
$ws.Cells.Item(205, 3).Value = "铁丝网栅栏"
$ws.Cells.Item(205, 4).Value = "wire mesh fence"
$ws.Cells.Item(205, 5).Value = "金網のフェンス"

$ws.Cells.Item(206, 1).Value = "'203"
$ws.Cells.Item(206, 2).Value = "EA 23.238"
$ws.Cells.Item(206, 3).Value = "铁丝网栅栏"
$ws.Cells.Item(206, 4).Value = "wire mesh fence"
$ws.Cells.Item(206, 5).Value = "金網のフェンス"

$ws.Cells.Item(207, 1).Value = "'204"
$ws.Cells.Item(207, 2).Value = "EA 23.238"
$ws.Cells.Item(207, 3).Value = "石头栅栏"
$ws.Cells.Item(207, 4).Value = "stone fence"
$ws.Cells.Item(207, 5).Value = "石のフェンス"
